# Insert 3 new weekly price rows for "Femacal de La Calera" / Chirimoya
# right before the current first data block (old row 478), pushing the
# existing rows down by 3 (old 478:501 -> new 481:504).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above what is currently row 478.
$ws.Rows("478:480").Insert()

# Row 478 - Especial
$ws.Cells.Item(478, 1).Value  = 3
$ws.Cells.Item(478, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(478, 3).Value  = "Coquimbo"
$ws.Cells.Item(478, 4).Value  = 45267
$ws.Cells.Item(478, 5).Value  = 5
$ws.Cells.Item(478, 6).Value  = "Fruta"
$ws.Cells.Item(478, 7).Value  = 100107
$ws.Cells.Item(478, 8).Value  = "Otros"
$ws.Cells.Item(478, 9).Value  = 100107002
$ws.Cells.Item(478, 10).Value = "Chirimoya"
$ws.Cells.Item(478, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(478, 12).Value = "Especial"
$ws.Cells.Item(478, 13).Value = 50
$ws.Cells.Item(478, 14).Value = 22000
$ws.Cells.Item(478, 15).Value = 22000
$ws.Cells.Item(478, 16).Value = 22000
$ws.Cells.Item(478, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(478, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(478, 19).Value = 2200
$ws.Cells.Item(478, 20).Value = 10

# Row 479 - Primera
$ws.Cells.Item(479, 1).Value  = 3
$ws.Cells.Item(479, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(479, 3).Value  = "Coquimbo"
$ws.Cells.Item(479, 4).Value  = 45267
$ws.Cells.Item(479, 5).Value  = 5
$ws.Cells.Item(479, 6).Value  = "Fruta"
$ws.Cells.Item(479, 7).Value  = 100107
$ws.Cells.Item(479, 8).Value  = "Otros"
$ws.Cells.Item(479, 9).Value  = 100107002
$ws.Cells.Item(479, 10).Value = "Chirimoya"
$ws.Cells.Item(479, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(479, 12).Value = "Primera"
$ws.Cells.Item(479, 13).Value = 50
$ws.Cells.Item(479, 14).Value = 19000
$ws.Cells.Item(479, 15).Value = 19000
$ws.Cells.Item(479, 16).Value = 19000
$ws.Cells.Item(479, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(479, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(479, 19).Value = 1900
$ws.Cells.Item(479, 20).Value = 10

# Row 480 - Segunda
$ws.Cells.Item(480, 1).Value  = 3
$ws.Cells.Item(480, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(480, 3).Value  = "Coquimbo"
$ws.Cells.Item(480, 4).Value  = 45267
$ws.Cells.Item(480, 5).Value  = 5
$ws.Cells.Item(480, 6).Value  = "Fruta"
$ws.Cells.Item(480, 7).Value  = 100107
$ws.Cells.Item(480, 8).Value  = "Otros"
$ws.Cells.Item(480, 9).Value  = 100107002
$ws.Cells.Item(480, 10).Value = "Chirimoya"
$ws.Cells.Item(480, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(480, 12).Value = "Segunda"
$ws.Cells.Item(480, 13).Value = 48
$ws.Cells.Item(480, 14).Value = 16000
$ws.Cells.Item(480, 15).Value = 16000
$ws.Cells.Item(480, 16).Value = 16000
$ws.Cells.Item(480, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(480, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(480, 19).Value = 1600
$ws.Cells.Item(480, 20).Value = 10
